# Fill in the crossword answers on Лист6 ("КОМПЬЮТЕР" crossword) so that
# the checker formulas on Лист7 evaluate to 1 for every answer cell and the
# total score (Лист7!L14) becomes 40, which in turn makes Лист6!B15 show
# "Молодец" instead of "Подумай ещё".

$wb  = $excel.ActiveWorkbook
$ws6 = $wb.Worksheets.Item("Лист6")
$ws7 = $wb.Worksheets.Item("Лист7")

# --- Word: ДИСКЕТА (horizontal, row 4, B4:H4) ---
$ws6.Range("B4").Value = "д"
$ws6.Range("C4").Value = "и"
$ws6.Range("D4").Value = "с"
$ws6.Range("E4").Value = "к"
$ws6.Range("F4").Value = "е"
$ws6.Range("G4").Value = "т"
$ws6.Range("H4").Value = "а"

# --- Word: КОМПЬЮТЕР (vertical, column E, E4:E12) ---
$ws6.Range("E4").Value  = "к"
$ws6.Range("E5").Value  = "о"
$ws6.Range("E6").Value  = "м"
$ws6.Range("E7").Value  = "п"
$ws6.Range("E8").Value  = "ь"
$ws6.Range("E9").Value  = "ю"
$ws6.Range("E10").Value = "т"
$ws6.Range("E11").Value = "е"
$ws6.Range("E12").Value = "р"

# --- Word: МОНИТОР (horizontal, row 6, E6:K6) ---
$ws6.Range("E6").Value = "м"
$ws6.Range("F6").Value = "о"
$ws6.Range("G6").Value = "н"
$ws6.Range("H6").Value = "и"
$ws6.Range("I6").Value = "т"
$ws6.Range("J6").Value = "о"
$ws6.Range("K6").Value = "р"

# --- Word: ПРОЦЕССОР (vertical, column J, J4:J12) ---
$ws6.Range("J4").Value  = "п"
$ws6.Range("J5").Value  = "р"
$ws6.Range("J6").Value  = "о"
$ws6.Range("J7").Value  = "ц"
$ws6.Range("J8").Value  = "е"
$ws6.Range("J9").Value  = "с"
$ws6.Range("J10").Value = "с"
$ws6.Range("J11").Value = "о"
$ws6.Range("J12").Value = "р"

# --- Word: МЫШЬ (horizontal, row 8, B8:E8) ---
$ws6.Range("B8").Value = "м"
$ws6.Range("D8").Value = "ш"
$ws6.Range("C8").Value = "ы"
$ws6.Range("E8").Value = "ь"

# --- Word: ДИСК (horizontal, row 9, H9:K9) ---
$ws6.Range("H9").Value = "д"
$ws6.Range("I9").Value = "и"
$ws6.Range("J9").Value = "с"
$ws6.Range("K9").Value = "к"

# --- Word: ПРИНТЕР (horizontal, row 12, D12:J12) ---
$ws6.Range("D12").Value = "п"
$ws6.Range("E12").Value = "р"
$ws6.Range("F12").Value = "и"
$ws6.Range("G12").Value = "н"
$ws6.Range("H12").Value = "т"
$ws6.Range("I12").Value = "е"
$ws6.Range("J12").Value = "р"

# Update the saved selection on Лист6 to match the author's final cursor position
$ws6.Range("L2").Select()
